$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D29").Value = "[Paper Review] UNet++와 U-Net3+ 논문 리뷰 & 의료 데이터 적용"
$ws.Range("E29").Value = "https://blog.promedius.ai/unet-unet3/"

$ws.Range("D36").Value = "History of Class Activation MAP (CAM)"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/354"
